# cs-en-us-048pct.xlsx weekly refresh: new crime data collected.
# Updates the "Volume/Number" + "Report Covering the Week" header text,
# and refreshes the Week-to-Date/28-Day/Year-to-Date/2-Year crime table
# (rows 14-29) with the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (rich text runs - only the trailing run's text
# content changes; re-apply the same font so the run keeps rendering
# identically to its neighbours).
# ---------------------------------------------------------------------
$volRun = $ws.Range("A8").Characters(21, 2)
$volRun.Text = "34"
$volRun.Font.Name = "Andale WT"
$volRun.Font.Size = 10

$date1 = $ws.Range("C9").Characters(27, 9)
$date1.Text = "8/21/2023"
$date1.Font.Name = "Andale WT"
$date1.Font.Size = 10

$date2 = $ws.Range("C9").Characters(47, 9)
$date2.Text = "8/27/2023"
$date2.Font.Name = "Andale WT"
$date2.Font.Size = 10

# ---------------------------------------------------------------------
# Crime-stat table numeric updates (rows 14-29)
# ---------------------------------------------------------------------
$ws.Range("F14").Value = 2
$ws.Range("I14").Value = 12
$ws.Range("J14").Value = 12
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 33.333333333333
$ws.Range("M14").Value = 1100
$ws.Range("N14").Value = -52

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 22
$ws.Range("J15").Value = 21
$ws.Range("K15").Value = 4.761904761904
$ws.Range("L15").Value = -8.333333333333
$ws.Range("M15").Value = 57.142857142857
$ws.Range("N15").Value = -43.589743589743

$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 38
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = 22.580645161290
$ws.Range("I16").Value = 264
$ws.Range("J16").Value = 299
$ws.Range("K16").Value = -11.705685618729
$ws.Range("L16").Value = 30.049261083743
$ws.Range("M16").Value = 9.543568464730
$ws.Range("N16").Value = -64.130434782608

$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 64
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = 48.837209302325
$ws.Range("I17").Value = 450
$ws.Range("J17").Value = 417
$ws.Range("K17").Value = 7.913669064748
$ws.Range("L17").Value = 47.058823529411
$ws.Range("M17").Value = 82.186234817813
$ws.Range("N17").Value = -5.063291139240

$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 75
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 160
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 167
$ws.Range("K18").Value = -4.191616766467
$ws.Range("L18").Value = 23.076923076923
$ws.Range("M18").Value = 1.265822784810
$ws.Range("N18").Value = -81.693363844393

$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 57.142857142857
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 62.5
$ws.Range("I19").Value = 324
$ws.Range("J19").Value = 299
$ws.Range("K19").Value = 8.361204013377
$ws.Range("L19").Value = 10.958904109589
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 6.229508196721

$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -10.344827586206
$ws.Range("I20").Value = 226
$ws.Range("J20").Value = 191
$ws.Range("K20").Value = 18.324607329842
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 189.743589743590
$ws.Range("N20").Value = -38.419618528610

$ws.Range("C21").Value = 42
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 16.666666666666
$ws.Range("F21").Value = 210
$ws.Range("G21").Value = 149
$ws.Range("H21").Value = 40.939597315436
$ws.Range("I21").Value = 1458
$ws.Range("J21").Value = 1406
$ws.Range("K21").Value = 3.698435277382
$ws.Range("L21").Value = 35.376044568245
$ws.Range("M21").Value = 52.670157068062
$ws.Range("N21").Value = -48.297872340425

$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50

$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -60.606060606060
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 97
$ws.Range("H24").Value = -32.989690721649
$ws.Range("I24").Value = 675
$ws.Range("J24").Value = 801
$ws.Range("K24").Value = -15.730337078651
$ws.Range("L24").Value = 24.080882352941
$ws.Range("M24").Value = 1.656626506024

$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 58.333333333333
$ws.Range("F25").Value = 74
$ws.Range("G25").Value = 65
$ws.Range("H25").Value = 13.846153846153
$ws.Range("I25").Value = 589
$ws.Range("J25").Value = 593
$ws.Range("K25").Value = -0.674536256323
$ws.Range("L25").Value = 24.261603375527
$ws.Range("M25").Value = -9.800918836140

$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -71.428571428571
$ws.Range("I26").Value = 35
$ws.Range("J26").Value = 42
$ws.Range("K26").Value = -16.666666666666
$ws.Range("L26").Value = 6.060606060606

$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -28.571428571428
$ws.Range("I27").Value = 51
$ws.Range("J27").Value = 55
$ws.Range("K27").Value = -7.272727272727
$ws.Range("L27").Value = 54.545454545454

$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 26
$ws.Range("J28").Value = 42
$ws.Range("K28").Value = -38.095238095238
$ws.Range("L28").Value = -36.585365853658
$ws.Range("M28").Value = -3.703703703703
$ws.Range("N28").Value = -69.767441860465

$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -50
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 19
$ws.Range("J29").Value = 36
$ws.Range("K29").Value = -47.222222222222
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = -17.391304347826
$ws.Range("N29").Value = -74.666666666666

# ---------------------------------------------------------------------
# Cells whose figures collapsed to "0" / "***.*" this week: these two
# strings are shared strings already used elsewhere in the sheet
# (e.g. D14/E14), so we set them as text (forcing NumberFormat to
# text first so "0" isn't silently reinterpreted as numeric 0), then
# copy the number-format from a cell that already has the correct
# "text" style so the cell's style index matches its siblings.
# ---------------------------------------------------------------------
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("G14").PasteSpecial(-4122)

$ws.Range("H14").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("H14").PasteSpecial(-4122)

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("M27").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("M27").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "***.*"
$ws.Range("M27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$excel.CutCopyMode = $false
